$p = $ppt.ActivePresentation

# --- Slide 12 ("Model1"): TextBox 50 (shape index 4) ---
# ".  " -> "." and append a new paragraph "Accuracy(0.6171)   " (3 trailing spaces)
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(4)
$tr12 = $shp12.TextFrame.TextRange
$len12 = $tr12.Text.Length
$tr12.Characters($len12 - 1, 2).Text = ""
$shp12.TextFrame.TextRange.InsertAfter("`rAccuracy(0.6171)   ")

# --- Slide 13 ("Model2"): TextBox 50 (shape index 4) ---
# keep ".  " as-is, append "Accuracy(0.6171)   " (3 trailing spaces) paragraph,
# then append one more empty paragraph (matching endParaRPr formatting)
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(4)
$tr13 = $shp13.TextFrame.TextRange
$tr13.InsertAfter("`rAccuracy(0.6171)   `r")

# --- Slide 14 ("Model3"): TextBox 50 (shape index 4) ---
# ".  " -> "." and append a new paragraph "Accuracy(0.6171)     " (5 trailing spaces)
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(4)
$tr14 = $shp14.TextFrame.TextRange
$len14 = $tr14.Text.Length
$tr14.Characters($len14 - 1, 2).Text = ""
$shp14.TextFrame.TextRange.InsertAfter("`rAccuracy(0.6171)     ")
